# rial_cumulative.xlsx - "update database and change read_price algorithm"
#
# The oldest cumulative-period column ("6 ماهه منتهی به 1399/06", column D)
# is retired and a new trailing column ("12 ماهه منتهی به 1401/12") is
# appended at the end (column M). Every other period keeps its data but
# slides one column to the left. In addition, the previously-published
# figures for "12 ماهه منتهی به 1400/12" (now column I) are restated with
# updated numbers, matching a later republish of that report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the oldest period column (D). Excel slides E:M -> D:L,
#        carrying values, formats and column widths with it. ---
$ws.Columns.Item(4).Delete()

# --- 2. Open up a fresh column M for the new trailing period, cloning
#        the formatting (fills/borders/alignment) from the column that is
#        now its left neighbour (L). ---
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns.Item(13).ColumnWidth = 28.17

# --- 3. Header / publish-date labels for the new column ---
$ws.Cells.Item(8, 13).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value = "1402-02-23"

# --- 4. Restate the republished "12 ماهه منتهی به 1400/12" figures
#        (now column I after the shift). ---
$ws.Cells.Item(9, 9).Value = "1402-02-23 (8)"
$ws.Cells.Item(19, 9).Value = 1182099
$ws.Cells.Item(20, 9).Value = 2848155
$ws.Cells.Item(22, 9).Value = 2719994
$ws.Cells.Item(24, 9).Value = 2719994
$ws.Cells.Item(25, 9).Value = 1360
$ws.Cells.Item(27, 9).Value = 680

# --- 5. Fill in the new trailing period's figures (column M) ---
$ws.Cells.Item(11, 13).Value = 37692415
$ws.Cells.Item(12, 13).Value = -26317310
$ws.Cells.Item(13, 13).Value = 11375105
$ws.Cells.Item(14, 13).Value = -3335380
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(16, 13).Value = -42234
$ws.Cells.Item(17, 13).Value = 7997491
$ws.Cells.Item(18, 13).Value = -2491810
$ws.Cells.Item(19, 13).Value = 1152674
$ws.Cells.Item(20, 13).Value = 6658355
$ws.Cells.Item(21, 13).Value = -52040
$ws.Cells.Item(22, 13).Value = 6606315
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(24, 13).Value = 6606315
$ws.Cells.Item(25, 13).Value = 1652
$ws.Cells.Item(26, 13).Value = 4000000
$ws.Cells.Item(27, 13).Value = 1652
